$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '242.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.23%'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.60%'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.109'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.57%'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05643'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.99%'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.493'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.26%'
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8252'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.23%'
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8683'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.32%'
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1332'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.27%'
# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06923'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.05%'
# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02859'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.12%'
# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09369'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.16%'
# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001514'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.73%'
# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'CoinExToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04173'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-9.11%'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006085'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.97%'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.522'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.05%'
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.022'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.42%'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.215'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.45%'
# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'One'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0006010'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.56%'
# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3149'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.19%'
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.03245'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.59%'
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.06%'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.611'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.42%'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1374'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.06%'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001210'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.85%'
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004443'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-1.89%'
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001180'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '22.95%'
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001404'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.63%'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03704'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.62%'
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005790'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '68.47%'
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-22.72%'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002312'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-11.06%'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009560'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.22%'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005087'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.95%'
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.06%'
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1200'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '10.16%'
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002385'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-6.28%'
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.06%'
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.06%'
